$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 data, following the same shape as existing rows (e.g. row 9)
# Numeric cells
$ws.Range("A10").Value = 111906849
$ws.Range("B10").Value = 88967
$ws.Range("E10").Value = 6039940
$ws.Range("Q10").Value = 813178.8074009671
$ws.Range("R10").Value = 7316199.822832054
$ws.Range("S10").Value = 20

# Text cells
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").Value = "DD"
$ws.Range("F10").Value = "Mandarinfingersvamp"
$ws.Range("G10").Value = "Ramaria tridentina"
$ws.Range("H10").Value = "Schild"
$ws.Range("I10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("P10").Value = "Torrkölen (Torrkölen), Nb"
$ws.Range("T10").Value = "Norrbotten"
$ws.Range("U10").Value = "Boden"
$ws.Range("V10").Value = "Norrbotten"
$ws.Range("W10").Value = "Överluleå"

# Date-like text cells - force text format so Excel doesn't coerce to a date serial
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "2023-09-05"
$ws.Range("Y10").Style = "Normal"

$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = "2023-09-05"
$ws.Range("AA10").Style = "Normal"

$ws.Range("Z10").Value = "11:11"
$ws.Range("AB10").Value = "11:11"

# Boolean cells
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false

# Remaining text cells
$ws.Range("AT10").Value = ""
$ws.Range("AW10").Value = "Linnea Åsedahl"
$ws.Range("AX10").Value = "Linnea Åsedahl"
$ws.Range("AY10").Value = ""
